$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("B2").Formula = "=LOG10(17.3)"
$ws.Range("C2").Formula = "=LOG10(0.329/2.86)"
$ws.Range("D2").Formula = "=LOG10(0.0476/2.86)"
$ws.Range("E2").Formula = "=LOG10(1.24/2.86)"

# Row 3 updates
$ws.Range("B3").Formula = "=LOG10(17.3)"
$ws.Range("C3").Formula = "=LOG10(0.345/3.06)"
$ws.Range("D3").Formula = "=LOG10(0.0457/3.09)"
$ws.Range("E3").Formula = "=LOG10(1.23/3.09)"

# Row 5 updates
$ws.Range("B5").Formula = "=LOG10(24.8)"
$ws.Range("C5").Formula = "=LOG10(0.221/2.86)"
$ws.Range("D5").Formula = "=LOG10(0.0735/2.86)"
$ws.Range("E5").Formula = "=LOG10(1.62/2.86)"

# Clear B7:E10 (only the A column values remain)
$ws.Range("B7:E10").ClearContents()

# Update selection to E5
$ws.Range("E5").Select()
